$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020416295938007
$ws.Range("D2").Value = 1.026373181646233
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.031622550597911
$ws.Range("I2").Value = 1.030123242816486
$ws.Range("J2").Value = 1.025613529653269
$ws.Range("K2").Value = 1.029196018135634
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.034430135674624
$ws.Range("N2").Value = 1.027070018504895

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021400348210294
$ws.Range("D3").Value = 1.027113678471736
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.032833526293307
$ws.Range("I3").Value = 1.030336961472128
$ws.Range("J3").Value = 1.02623457961382
$ws.Range("K3").Value = 1.029744197124219
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.035448635651216
$ws.Range("N3").Value = 1.027691950427623

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022036788193595
$ws.Range("D4").Value = 1.027592165882076
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.033616829310957
$ws.Range("I4").Value = 1.030473249913109
$ws.Range("J4").Value = 1.026635568282763
$ws.Range("K4").Value = 1.030097580749151
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.036106835780602
$ws.Range("N4").Value = 1.02809350854647

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022304274051802
$ws.Range("D5").Value = 1.027793162494624
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.033946064049698
$ws.Range("I5").Value = 1.030530066034013
$ws.Range("J5").Value = 1.026803934940812
$ws.Range("K5").Value = 1.030245825664211
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.036383343335093
$ws.Range("N5").Value = 1.028262114304485

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022349181797588
$ws.Range("D6").Value = 1.027826901338637
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.034001340178454
$ws.Range("I6").Value = 1.030539577581413
$ws.Range("J6").Value = 1.02683219218369
$ws.Range("K6").Value = 1.030270698004229
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.036429758475459
$ws.Range("N6").Value = 1.028290411675889

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022040362641411
$ws.Range("D7").Value = 1.027594852236186
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.03362122881917
$ws.Range("I7").Value = 1.030474010977554
$ws.Range("J7").Value = 1.02663781882767
$ws.Range("K7").Value = 1.030099562853277
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.036110531272218
$ws.Range("N7").Value = 1.028095762287409

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020748924795312
$ws.Range("D8").Value = 1.026623573027708
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.032031863925194
$ws.Range("I8").Value = 1.03019588434908
$ws.Range("J8").Value = 1.025823596883825
$ws.Range("K8").Value = 1.0293815518751
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.034774516342468
$ws.Range("N8").Value = 1.027280384055011

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.0184708829989
$ws.Range("D9").Value = 1.024907000218624
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.029229015656958
$ws.Range("I9").Value = 1.029690470761215
$ws.Range("J9").Value = 1.024382160812531
$ws.Range("K9").Value = 1.02810619455814
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.032413840623586
$ws.Range("N9").Value = 1.02583690097916

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.016950585792647
$ws.Range("D10").Value = 1.023759253726228
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.027358913980313
$ws.Range("I10").Value = 1.029343244027463
$ws.Range("J10").Value = 1.023416723387171
$ws.Range("K10").Value = 1.027249164808804
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.030835672921002
$ws.Range("N10").Value = 1.02487009252192

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016291895986903
$ws.Range("D11").Value = 1.023261475029037
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.026548758755526
$ws.Range("I11").Value = 1.029190454348271
$ws.Range("J11").Value = 1.022997616177442
$ws.Range("K11").Value = 1.026876453847202
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.030151256936482
$ws.Range("N11").Value = 1.02445039013188

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016047170019769
$ws.Range("D12").Value = 1.023076458533907
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.026247771084044
$ws.Range("I12").Value = 1.029133335275024
$ws.Range("J12").Value = 1.022841780807762
$ws.Range("K12").Value = 1.026737770458211
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.029896874090527
$ws.Range("N12").Value = 1.024294333458101

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01609966725937
$ws.Range("D13").Value = 1.023116150593501
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.026312336698756
$ws.Range("I13").Value = 1.029145604078687
$ws.Range("J13").Value = 1.022875215278313
$ws.Range("K13").Value = 1.026767529461964
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.029951447346376
$ws.Range("N13").Value = 1.024327815409436

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016271668073557
$ws.Range("D14").Value = 1.023246183939723
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.026523880252519
$ws.Range("I14").Value = 1.029185740337162
$ws.Range("J14").Value = 1.022984738047509
$ws.Range("K14").Value = 1.026864995169707
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.030130232847518
$ws.Range("N14").Value = 1.024437493713526

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0163776355866
$ws.Range("D15").Value = 1.023326285976234
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.026654211228135
$ws.Range("I15").Value = 1.029210421091588
$ws.Range("J15").Value = 1.023052197397883
$ws.Range("K15").Value = 1.026925014953706
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.030240367221882
$ws.Range("N15").Value = 1.024505048863915

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.016994292584457
$ws.Range("D16").Value = 1.023792272903852
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.027412672991497
$ws.Range("I16").Value = 1.02935333281624
$ws.Range("J16").Value = 1.023444515677939
$ws.Range("K16").Value = 1.027273866470481
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.030881072944418
$ws.Range("N16").Value = 1.024897924280929

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017381000180749
$ws.Range("D17").Value = 1.024084361128731
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.027888330834929
$ws.Range("I17").Value = 1.029442324881328
$ws.Range("J17").Value = 1.023690321007975
$ws.Range("K17").Value = 1.027492260499639
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.03128268637582
$ws.Range("N17").Value = 1.025144078682728

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01760652247856
$ws.Range("D18").Value = 1.024254654333901
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.028165736591434
$ws.Range("I18").Value = 1.02949399718683
$ws.Range("J18").Value = 1.023833592208224
$ws.Range("K18").Value = 1.027619490553001
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.031516838601541
$ws.Range("N18").Value = 1.025287553344514

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017683413374185
$ws.Range("D19").Value = 1.024312706857649
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.028260318485469
$ws.Range("I19").Value = 1.029511576199345
$ws.Range("J19").Value = 1.023882426535153
$ws.Range("K19").Value = 1.027662846334767
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.031596661178274
$ws.Range("N19").Value = 1.025336457021788

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017339514018322
$ws.Range("D20").Value = 1.02405303077911
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.027837301130768
$ws.Range("I20").Value = 1.029432801208018
$ws.Range("J20").Value = 1.023663959065125
$ws.Range("K20").Value = 1.027468844968864
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.031239607623787
$ws.Range("N20").Value = 1.025117679302895

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016221019760209
$ws.Range("D21").Value = 1.02320789565101
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.026461587615721
$ws.Range("I21").Value = 1.029173931316854
$ws.Range("J21").Value = 1.022952490747205
$ws.Range("K21").Value = 1.026836300625658
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.030077589407463
$ws.Range("N21").Value = 1.024405200618356

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015517435005414
$ws.Range("D22").Value = 1.022675835473592
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.025596274317363
$ws.Range("I22").Value = 1.029009050976622
$ws.Range("J22").Value = 1.022504234940262
$ws.Range("K22").Value = 1.026437195153753
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.029346054536984
$ws.Range("N22").Value = 1.023956308236753

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015890451276247
$ws.Range("D23").Value = 1.022957955931129
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.026055026731559
$ws.Range("I23").Value = 1.029096657938545
$ws.Range("J23").Value = 1.022741951761423
$ws.Range("K23").Value = 1.026648901145663
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.029733943363829
$ws.Range("N23").Value = 1.024194362643066

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017358259945609
$ws.Range("D24").Value = 1.024067187853078
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.027860359372645
$ws.Range("I24").Value = 1.029437105272888
$ws.Range("J24").Value = 1.023675871210131
$ws.Range("K24").Value = 1.027479425919637
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.031259073372561
$ws.Range("N24").Value = 1.025129608364513

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019060092851144
$ws.Range("D25").Value = 1.025351370905135
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.029953884881316
$ws.Range("I25").Value = 1.029822946059379
$ws.Range("J25").Value = 1.024755596665057
$ws.Range("K25").Value = 1.028437103430592
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.033024900902748
$ws.Range("N25").Value = 1.026210867153432
